$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force text number format so Excel does not auto-coerce numeric-looking
    # strings (e.g. "1.007") into floating point numbers on assignment.
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "24.848.06"
$ws.Range("E2").Value = "  +2.10%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.704.08"
$ws.Range("E3").Value = "  +1.82%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.007"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "312.03"
$ws.Range("E5").Value = "  +2.18%  "

# Row 6 - USDC
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  +0.02%  "

# Row 7 - XRP
Set-TextValue "D7" "0.3729"
$ws.Range("E7").Value = "  +1.38%  "

# Row 8 - OKB
Set-TextValue "D8" "49.27"
$ws.Range("E8").Value = "  +3.64%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.3425"
$ws.Range("E9").Value = "  +0.30%  "

# Row 10 - Polygon
Set-TextValue "D10" "1.217"
$ws.Range("E10").Value = "  +5.47%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.07498"
$ws.Range("E11").Value = "  +4.27%  "

# Row 12 - BinanceUSD
Set-TextValue "D12" "1.002"
$ws.Range("E12").Value = "  -0.08%  "

# Row 13 - Solana
Set-TextValue "D13" "21.06"
$ws.Range("E13").Value = "  +5.39%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.315"
$ws.Range("E14").Value = "  +3.30%  "

# Row 15 - Chainlink
Set-TextValue "D15" "7.028"
$ws.Range("E15").Value = "  +4.94%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "1.707.61"
$ws.Range("E16").Value = "  +1.96%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.00001128"
$ws.Range("E17").Value = "  +2.82%  "

# Row 18 - TRON
Set-TextValue "D18" "0.06717"
$ws.Range("E18").Value = "  +1.13%  "

# Row 19 - Dai
Set-TextValue "D19" "0.9995"
$ws.Range("E19").Value = "  -0.04%  "

# Row 20 - Litecoin
Set-TextValue "D20" "83.77"
$ws.Range("E20").Value = "  +4.82%  "

# Row 21 - Avalanche
Set-TextValue "D21" "17.26"
$ws.Range("E21").Value = "  +5.30%  "

# Row 22 - Uniswap
Set-TextValue "D22" "6.351"
$ws.Range("E22").Value = "  +4.58%  "

# Row 23 - Cosmos
Set-TextValue "D23" "13.10"
$ws.Range("E23").Value = "  +7.76%  "

# Row 24 - WrappedBTC
Set-TextValue "D24" "24.863.78"
$ws.Range("E24").Value = "  +2.34%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.457"
$ws.Range("E25").Value = "  +0.84%  "

# Row 26 - LidoDAOToken
Set-TextValue "D26" "2.775"
$ws.Range("E26").Value = "  +5.48%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "20.29"
$ws.Range("E27").Value = "  +4.95%  "

# Row 28 - Monero
Set-TextValue "D28" "149.48"
$ws.Range("E28").Value = "  -1.26%  "

# Row 29 - BitcoinCash
Set-TextValue "D29" "132.31"
$ws.Range("E29").Value = "  +3.82%  "

# Row 30 - ImmutableX
Set-TextValue "D30" "1.255"
$ws.Range("E30").Value = "  +30.34%  "

# Row 31 - WrappedliquidstakedEther2.0
Set-TextValue "D31" "1.894.41"
$ws.Range("E31").Value = "  +1.85%  "

# Row 32 - Filecoin
Set-TextValue "D32" "6.764"
$ws.Range("E32").Value = "  +8.50%  "

# Row 33 - HuobiToken
Set-TextValue "D33" "4.231"
$ws.Range("E33").Value = "  +4.49%  "

# Row 34 - Aptos
Set-TextValue "D34" "13.72"
$ws.Range("E34").Value = "  +12.12%  "

# Row 35 - Stellar
Set-TextValue "D35" "0.08760"
$ws.Range("E35").Value = "  +3.81%  "

# Row 36 - WEMIXTOKEN
Set-TextValue "D36" "1.770"
$ws.Range("E36").Value = "  +5.62%  "

# Row 37 - InternetComputer(DFINITY)
Set-TextValue "D37" "5.583"
$ws.Range("E37").Value = "  +6.16%  "

# Row 38 - Hedera
Set-TextValue "D38" "0.06646"
$ws.Range("E38").Value = "  +4.46%  "

# Row 39 - FraxShare
Set-TextValue "D39" "9.119"
$ws.Range("E39").Value = "  +5.65%  "

# Row 40 - was Algorand, now VeChain (rows 40/41 swap coin order)
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.02408"
$ws.Range("E40").Value = "  +4.43%  "

# Row 41 - was VeChain, now Algorand
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D41" "0.2236"
$ws.Range("E41").Value = "  +7.80%  "

# Row 42 - TrustWalletToken
$ws.Range("E42").Value = "  +3.12%  "

# Row 43 - TheSandbox
Set-TextValue "D43" "0.6419"
$ws.Range("E43").Value = "  +6.20%  "

# Row 44 - Frax
Set-TextValue "D44" "1.000"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "13.86"
$ws.Range("E45").Value = "  +7.14%  "

# Row 46 - Decentraland
Set-TextValue "D46" "0.6130"
$ws.Range("E46").Value = "  +4.86%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "3.835"
$ws.Range("E47").Value = "  +2.56%  "

# Row 48 - NEARProtocol
Set-TextValue "D48" "2.110"
$ws.Range("E48").Value = "  +5.27%  "

# Row 49 - Quant
Set-TextValue "D49" "129.68"
$ws.Range("E49").Value = "  +3.45%  "

# Row 50 - Cronos
Set-TextValue "D50" "0.07305"
$ws.Range("E50").Value = "  +2.44%  "

# Row 51 - Aave
Set-TextValue "D51" "79.71"
$ws.Range("E51").Value = "  +5.59%  "
